$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$rng = $ws.Range("A1")

# Step 1: replace entire cell text with updated content (this resets rich-text runs to plain)
$rng.Value2 = "DO NOT DELETE THIS ROW! RETAIN THE HEADING ROW!`nInstructions: Starting on Row 3, fill in the relevant fields. Do not delete rows 1 and 2.`nFor SPOUSE EMAIL and SPOUSE ID, these are mutually exclusive. If using them (they're not mandatory), use one or the other but not both. If both are used, the ID will take precendence and the email will be ignored.`n`nNAME: Full Name`nEMAIL: Propely formatted email address. NOTE, IF AN EMAIL ADDRESS ALREADY EXISTS FOR A USER IN THE SYSTEM, IT WILL UPDATE THAT USER RECORD.`nMOBILE PHONE: In the format 04XXXXXXXX (spaces can be used)`nGENDER: male, female, m or f`nYEAR OF BIRTH: Optional. 4 digit year. Eg: 1985`nAPPOINTMENT: Optional. Allowed values only: elder, ministerial servant`nSERVING AS: Optional. Allowed values only: field missionary, special pioneer, bethel family member, regular pioneer, publisher`nMARITAL STATUS: Optional. Allowed values only: single, married, separated, divorced, widowed`nSPOUSE EMAIL: Optional. Used to link spouses together. If a matching email is found, it will attach the users`nSPOUSE ID: Optional. Used to link a user that already exists in the system to this user`nRESPONSIBLE BROTHER: Inidcates in the system that a user (brother) has been trained to oversee a shift. Allowed values only. TRUE, FALSE.`nIS UNRESTRICTED: TRUE is the default. If set to false (i.e. indicating they're a 'restricted' user), the volunteer cannot self-roster and they cannot see any shifts other than those they've been rostered onto. Allowed values only. TRUE, FALSE."

# Step 2: reapply per-run character formatting to recreate the rich text runs
$r = $rng.Characters(1, 352)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Font.Underline = $false

$r = $rng.Characters(353, 1)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $false
$r.Font.Underline = $false

$r = $rng.Characters(354, 5)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Font.Underline = $false

$r = $rng.Characters(359, 11)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $false
$r.Font.Underline = $false

$r = $rng.Characters(370, 40)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Font.Underline = $false

$r = $rng.Characters(410, 98)
$r.Font.Name = "Calibri (Body)"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Font.Underline = $true

$r = $rng.Characters(508, 1)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Font.Underline = $false

$r = $rng.Characters(509, 1)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $false
$r.Font.Underline = $false

$r = $rng.Characters(510, 13)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Font.Underline = $false

$r = $rng.Characters(523, 47)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $false
$r.Font.Underline = $false

$r = $rng.Characters(570, 7)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Font.Underline = $false

$r = $rng.Characters(577, 22)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $false
$r.Font.Underline = $false

$r = $rng.Characters(599, 14)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Font.Underline = $false

$r = $rng.Characters(613, 34)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $false
$r.Font.Underline = $false

$r = $rng.Characters(647, 12)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Font.Underline = $false

$r = $rng.Characters(659, 59)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $false
$r.Font.Underline = $false

$r = $rng.Characters(718, 11)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Font.Underline = $false

$r = $rng.Characters(729, 116)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $false
$r.Font.Underline = $false

$r = $rng.Characters(845, 15)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Font.Underline = $false

$r = $rng.Characters(860, 78)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $false
$r.Font.Underline = $false

$r = $rng.Characters(938, 13)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Font.Underline = $false

$r = $rng.Characters(951, 97)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $false
$r.Font.Underline = $false

$r = $rng.Characters(1048, 10)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Font.Underline = $false

$r = $rng.Characters(1058, 78)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $false
$r.Font.Underline = $false

$r = $rng.Characters(1136, 20)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Font.Underline = $false

$r = $rng.Characters(1156, 118)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $false
$r.Font.Underline = $false

$r = $rng.Characters(1274, 16)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $true
$r.Font.Underline = $false

$r = $rng.Characters(1290, 227)
$r.Font.Name = "Calibri"
$r.Font.Size = 12
$r.Font.Bold = $false
$r.Font.Underline = $false
